$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 = theta_se values
$ws.Range("B4").Value = "(0.25)"
$ws.Range("C4").Value = "(0.37)"
$ws.Range("D4").Value = "(0.01)"
$ws.Range("E4").Value = "(0.22)"
$ws.Range("F4").Value = "(0.05)"

# Row 6 = lambda_se values
$ws.Range("B6").Value = "(0.22)"
$ws.Range("C6").Value = "(0.12)"
$ws.Range("D6").Value = "(0.2)"
$ws.Range("E6").Value = "(0.22)"
$ws.Range("F6").Value = "(0.05)"

# Row 7 = total_dof, multiple_imputation column count changed
$ws.Range("F7").Value = 6742
